$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Лист1" -> "Общий список"
$ws.Name = "Общий список"

# Append new document rows under the existing list
# (write order chosen to match the shared-string table build order)
$ws.Range("A17").Value = "UML Class Diagram"
$ws.Range("A16").Value = "UML Use Case Diagram"
$ws.Range("A18").Value = "Class Analysis Diagram"
$ws.Range("A19").Value = "IDEF0"

# Move selection to A20 to match the post-edit cursor position
$ws.Range("A20").Select()
